# "Generate Report for Handback"
#
# The handback process discovered that the transformed handback file name
# (ucirm0kp.h0w) didn't match the expected handoff-derived name for the
# e09a1937-... source file, in both the zh-cn and de-de targets. Update the
# status report workbook accordingly:
#   - Flip that row's Status from "Ready for handoff" to
#     "Handback transform failed" (Overview + both per-language sheets).
#   - Populate the (previously empty) "Error Detail" column for that row on
#     both per-language sheets with the mismatch diagnostic.
#   - Widen the "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

$zhcn = "Handback file name: ucirm0kp.h0w is different with handoff file name: e09a1937-ee52-4a29-8d08-9b8019e5dbaa.3d2c52763c1c92ccf14cd1e53b6b61f4937f73b5.zh-cn."
$dede = "Handback file name: ucirm0kp.h0w is different with handoff file name: e09a1937-ee52-4a29-8d08-9b8019e5dbaa.3d2c52763c1c92ccf14cd1e53b6b61f4937f73b5.de-de."

# --- Overview sheet: update the per-language status cells for the
#     e09a1937 row (E3 = zh-cn status, F3 = de-de status) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"

# --- zh-cn sheet: Status (C3) + Error Detail (P3) + column width ---
# (ColumnWidth 39.17 round-trips to a stored column width of exactly 40,
# matching the other full-width columns on this sheet, e.g. column A/G.)
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Handback transform failed"
$zh.Range("P3").Value = $zhcn
$zh.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: Status (C3) + Error Detail (P3) + column width ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Handback transform failed"
$de.Range("P3").Value = $dede
$de.Columns.Item(16).ColumnWidth = 39.17
